$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row: E1 rename, add BB1:BE1 new headers ---
$ws.Range("E1").Value = "YesterdaySalesQty"

# Copy header style (bold, border, centered) from an existing header cell to the new header cells
$ws.Range("BA1").Copy()
$ws.Range("BB1:BE1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$hdrArr = New-Object 'object[,]' 1,4
$hdrArr[0,0] = "TP"
$hdrArr[0,1] = "TP Sales Value"
$hdrArr[0,2] = "Net Sales Value"
$hdrArr[0,3] = "Discount"
$ws.Range("BB1:BE1").Value = $hdrArr

# --- Update data rows: B (brand), C (item), D (uom), and new BB:BE values ---
$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Desodin'
$bcd[0,1] = 'Desodin 60ml Syrup'
$bcd[0,2] = '60 ml'
$ws.Range("B2:D2").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 18.74
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB2:BE2").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Dinafex'
$bcd[0,1] = 'Dinafex 120mg Tablet'
$bcd[0,2] = '30''s'
$ws.Range("B3:D3").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 179.91
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB3:BE3").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Dinafex'
$bcd[0,1] = 'Dinafex 180mg Tablet'
$bcd[0,2] = '30''s'
$ws.Range("B4:D4").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 224.89
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB4:BE4").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Dinafex'
$bcd[0,1] = 'Dinafex 60mg Tablet'
$bcd[0,2] = '30''s'
$ws.Range("B5:D5").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 78.71
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB5:BE5").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Dorenta'
$bcd[0,1] = 'Dorenta 50mg Tablet'
$bcd[0,2] = '50''s'
$ws.Range("B6:D6").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 93.71
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB6:BE6").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Etorix'
$bcd[0,1] = 'Etorix 120mg Tablet'
$bcd[0,2] = '20''s'
$ws.Range("B7:D7").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 209.9
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB7:BE7").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Etorix'
$bcd[0,1] = 'Etorix 90mg Tablet'
$bcd[0,2] = '30''s'
$ws.Range("B8:D8").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 269.87
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB8:BE8").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Etorix'
$bcd[0,1] = 'Etorix 60mg Tablet - 40''s'
$bcd[0,2] = '40''s'
$ws.Range("B9:D9").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 209.9
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB9:BE9").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Fenobac'
$bcd[0,1] = 'Fenobac 100ml Syrup'
$bcd[0,2] = '100ml'
$ws.Range("B10:D10").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 74.96
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB10:BE10").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Flucloxin'
$bcd[0,1] = 'Flucloxin 500mg Capsule'
$bcd[0,2] = '30 ''s'
$ws.Range("B11:D11").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 237.74
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB11:BE11").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Flucloxin'
$bcd[0,1] = 'Flucloxin 500mg Capsule - 36''s'
$bcd[0,2] = '36 ''s'
$ws.Range("B12:D12").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 284.21
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB12:BE12").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Geminox'
$bcd[0,1] = 'Geminox 320mg Tablet - 8''s'
$bcd[0,2] = '8 ''s'
$ws.Range("B13:D13").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 389.8
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB13:BE13").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Ketonic'
$bcd[0,1] = 'Ketonic 10mg Tablet'
$bcd[0,2] = '20''s'
$ws.Range("B14:D14").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 150.38
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB14:BE14").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Ketonic'
$bcd[0,1] = 'Ketonic 30mg Injection'
$bcd[0,2] = '5 ''s'
$ws.Range("B15:D15").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 206.77
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB15:BE15").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Ketonic'
$bcd[0,1] = 'Ketonic 30mg IM/IV Injection - 4''s'
$bcd[0,2] = '4''s'
$ws.Range("B16:D16").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 165.41
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB16:BE16").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Kynol'
$bcd[0,1] = 'Kynol D 25mg Tablet'
$bcd[0,2] = '60 ''s'
$ws.Range("B17:D17").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 180.45
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB17:BE17").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Kynol'
$bcd[0,1] = 'Kynol TR 100mg Capsule'
$bcd[0,2] = '50 ''s'
$ws.Range("B18:D18").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 262.37
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB18:BE18").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Kynol'
$bcd[0,1] = 'Kynol TR 200mg Capsule'
$bcd[0,2] = '30 ''s'
$ws.Range("B19:D19").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 224.89
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB19:BE19").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Naprox'
$bcd[0,1] = 'Naprox Plus 500mg Tablet - 30''s'
$bcd[0,2] = '30 ''s'
$ws.Range("B20:D20").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 224.89
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB20:BE20").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Oradin'
$bcd[0,1] = 'Oradin Plus Tablet - 40''s'
$bcd[0,2] = '40 ''s'
$ws.Range("B21:D21").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 209.9
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB21:BE21").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Osticare'
$bcd[0,1] = 'Osticare Tablet 24''s'
$bcd[0,2] = '24''s'
$ws.Range("B22:D22").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 215.89
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB22:BE22").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Sk-Mox'
$bcd[0,1] = 'Sk-Mox 500mg Capsule'
$bcd[0,2] = '48 ''s'
$ws.Range("B23:D23").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 219.13
$bbbe[0,1] = 219.13
$bbbe[0,2] = 219.13
$bbbe[0,3] = 0
$ws.Range("BB23:BE23").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Zithrox'
$bcd[0,1] = 'Zithrox 500mg Tablet'
$bcd[0,2] = '6 ''s'
$ws.Range("B24:D24").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 136.83
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB24:BE24").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Zithrox'
$bcd[0,1] = 'Zithrox 15ml Suspension'
$bcd[0,2] = '15 ml'
$ws.Range("B25:D25").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 71.96
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB25:BE25").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Zithrox'
$bcd[0,1] = 'Zithrox 30ml Dry Suspension'
$bcd[0,2] = '30ml'
$ws.Range("B26:D26").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 97.45
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB26:BE26").Value = $bbbe

$bcd = New-Object 'object[,]' 1,3
$bcd[0,0] = 'Zithrox'
$bcd[0,1] = 'Zithrox 250mg Tablet - 6''s'
$bcd[0,2] = '6''s'
$ws.Range("B27:D27").Value = $bcd
$bbbe = New-Object 'object[,]' 1,4
$bbbe[0,0] = 89.96
$bbbe[0,1] = 0
$bbbe[0,2] = 0
$bbbe[0,3] = 0
$ws.Range("BB27:BE27").Value = $bbbe

